# Update countries & provincias Spain
# Applies the data refresh recorded in the commit "Update countries & provincias Spain":
#  - updates the "Datos actualizados a ..." timestamp
#  - refreshes COVID case counters for several countries
#  - swaps the Montserrat / Islas Malvinas rows (their totals tied, so the
#    relative order produced by the source data changed) and updates their
#    Recuperados/Muertes values

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Timestamp at the top of the sheet ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 11:24"

# --- Country rows with refreshed figures ---
# Row 25: Filipinas
$ws.Range("B25").Value = 164474
$ws.Range("C25").Value = 3314
$ws.Range("D25").Value = 112759
$ws.Range("E25").Value = 49034
$ws.Range("G25").Value = 18
$ws.Range("H25").Value = 2681

# Row 26: Indonesia
$ws.Range("B26").Value = 141370
$ws.Range("C26").Value = 1821
$ws.Range("D26").Value = 94458
$ws.Range("E26").Value = 40705
$ws.Range("G26").Value = 57
$ws.Range("H26").Value = 6207

# Row 47: Polonia
$ws.Range("B47").Value = 57279
$ws.Range("C47").Value = 595
$ws.Range("D47").Value = 39359
$ws.Range("E47").Value = 16035
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 1885

# Row 72: Austria
$ws.Range("B72").Value = 23534
$ws.Range("C72").Value = 164
$ws.Range("D72").Value = 20765
$ws.Range("E72").Value = 2040
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 729

# Row 90: Malasia
$ws.Range("B90").Value = 9212
$ws.Range("C90").Value = 12
$ws.Range("D90").Value = 8876
$ws.Range("E90").Value = 211

# Row 111: Hong Kong
$ws.Range("B111").Value = 4525
$ws.Range("C111").Value = 44
$ws.Range("D111").Value = 3599
$ws.Range("E111").Value = 857

# Row 123: Eslovaquia
$ws.Range("B123").Value = 2907
$ws.Range("C123").Value = 5
$ws.Range("E123").Value = 907

# Row 124: Sri Lanka
$ws.Range("D124").Value = 2676
$ws.Range("E124").Value = 206

# Row 168: Taiwan
$ws.Range("B168").Value = 485
$ws.Range("C168").Value = 1
$ws.Range("E168").Value = 28

# --- Montserrat / Islas Malvinas swap (rows 213-214) ---
# Row 213 becomes Islas Malvinas, row 214 becomes Montserrat.
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
